$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "datos actualizados" timestamp cell (A1): 19:05 -> 19:35
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 19:35"

# Update country rows: resorted order (name + stats) and/or refreshed stats
# Columns: A=Pais  B=Casos totales  C=Nuevos casos  D=Casos activos
#          E=Recuperados  F=Casos criticos  G=Muertes hoy  H=Muertes

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 1655670
$ws.Range("C4").Value = 10576
$ws.Range("D4").Value = 438562
$ws.Range("E4").Value = 1118963
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 498
$ws.Range("H4").Value = 98145

# Row 7: España
$ws.Range("A7").Value = "España"
$ws.Range("B7").Value = 281904
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 196958
$ws.Range("E7").Value = 56268
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 50
$ws.Range("H7").Value = 28678

# Row 25: Ecuador
$ws.Range("A25").Value = "Ecuador"
$ws.Range("B25").Value = 36258
$ws.Range("C25").Value = 430
$ws.Range("D25").Value = 3557
$ws.Range("E25").Value = 29605
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = 3096

# Row 30: Suiza
$ws.Range("A30").Value = "Suiza"
$ws.Range("B30").Value = 30725
$ws.Range("C30").Value = 18
$ws.Range("D30").Value = 28000
$ws.Range("E30").Value = 820
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 1905

# Row 42: Egipto
$ws.Range("A42").Value = "Egipto"
$ws.Range("B42").Value = 16513
$ws.Range("C42").Value = 727
$ws.Range("D42").Value = 4628
$ws.Range("E42").Value = 11150
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 28
$ws.Range("H42").Value = 735

# Row 43: Japon
$ws.Range("A43").Value = "Japon"
$ws.Range("B43").Value = 16513
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 13005
$ws.Range("E43").Value = 2712
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 796

# Row 44: Austria
$ws.Range("A44").Value = "Austria"
$ws.Range("B44").Value = 16486
$ws.Range("C44").Value = 50
$ws.Range("D44").Value = 15037
$ws.Range("E44").Value = 810
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 4
$ws.Range("H44").Value = 639

# Row 45: Republica Dominicana
$ws.Range("A45").Value = "Republica Dominicana"
$ws.Range("B45").Value = 14422
$ws.Range("C45").Value = 433
$ws.Range("D45").Value = 7854
$ws.Range("E45").Value = 6110
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 458

# Row 117: Paraguay
$ws.Range("A117").Value = "Paraguay"
$ws.Range("B117").Value = 850
$ws.Range("C117").Value = 12
$ws.Range("D117").Value = 298
$ws.Range("E117").Value = 541
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 11

# Row 126: Republica del Chad
$ws.Range("A126").Value = "Republica del Chad"
$ws.Range("B126").Value = 648
$ws.Range("C126").Value = 37
$ws.Range("D126").Value = 204
$ws.Range("E126").Value = 384
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 2
$ws.Range("H126").Value = 60

# Row 127: Sierra Leona
$ws.Range("A127").Value = "Sierra Leona"
$ws.Range("B127").Value = 621
$ws.Range("C127").Value = 15
$ws.Range("D127").Value = 241
$ws.Range("E127").Value = 341
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 39

# Row 144: Ruanda
$ws.Range("A144").Value = "Ruanda"
$ws.Range("B144").Value = 325
$ws.Range("C144").Value = 4
$ws.Range("D144").Value = 227
$ws.Range("E144").Value = 98
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 0

# Row 145: Vietnam
$ws.Range("A145").Value = "Vietnam"
$ws.Range("B145").Value = 324
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 267
$ws.Range("E145").Value = 57
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 0

# Row 146: Montenegro
$ws.Range("A146").Value = "Montenegro"
$ws.Range("B146").Value = 324
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 314
$ws.Range("E146").Value = 1
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 9

# Row 151: Suazilandia
$ws.Range("A151").Value = "Suazilandia"
$ws.Range("B151").Value = 238
$ws.Range("C151").Value = 13
$ws.Range("D151").Value = 119
$ws.Range("E151").Value = 117
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 2

# Row 181: Zimbabue
$ws.Range("A181").Value = "Zimbabue"
$ws.Range("B181").Value = 56
$ws.Range("C181").Value = 5
$ws.Range("D181").Value = 18
$ws.Range("E181").Value = 34
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 4
